$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.942.24'
$ws.Range("E2").Value = '  -0.62%  '
$ws.Range("D3").Value = '1.818.40'
$ws.Range("E3").Value = '  -0.51%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''230.20'
$ws.Range("E5").Value = '  -1.12%  '
$ws.Range("D6").Value = '''0.616'
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '''40.10'
$ws.Range("E8").Value = '  -6.75%  '
$ws.Range("E9").Value = '  +4.13%  '
$ws.Range("D10").Value = '''0.0682'
$ws.Range("E10").Value = '  -1.01%  '
$ws.Range("E11").Value = '  -1.46%  '
$ws.Range("D12").Value = '2.080.49'
$ws.Range("E12").Value = '  -0.63%  '
$ws.Range("D13").Value = '''11.32'
$ws.Range("E13").Value = '  +1.00%  '
$ws.Range("D14").Value = '''0.669'
$ws.Range("E14").Value = '  +0.54%  '
$ws.Range("D15").Value = '1.813.31'
$ws.Range("E15").Value = '  -0.87%  '
$ws.Range("E16").Value = '  -1.88%  '
$ws.Range("D17").Value = '34.954.47'
$ws.Range("E17").Value = '  -0.55%  '
$ws.Range("D18").Value = '''69.65'
$ws.Range("E18").Value = '  -0.64%  '
$ws.Range("D19").Value = '0.0₃0784'
$ws.Range("E19").Value = '  -1.05%  '
$ws.Range("D20").Value = '''240.68'
$ws.Range("E20").Value = '  +0.26%  '
$ws.Range("D21").Value = '''12.03'
$ws.Range("E21").Value = '  +1.55%  '
$ws.Range("D22").Value = '''4.66'
$ws.Range("E22").Value = '  +1.20%  '
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("E24").Value = '  +2.18%  '
$ws.Range("D25").Value = '''173.17'
$ws.Range("E25").Value = '  +0.94%  '
$ws.Range("D26").Value = '''7.83'
$ws.Range("E26").Value = '  +0.34%  '
$ws.Range("D27").Value = '''0.123'
$ws.Range("E27").Value = '  +2.11%  '
$ws.Range("D28").Value = '''17.34'
$ws.Range("E28").Value = '  -1.41%  '
$ws.Range("E29").Value = '  -5.67%  '
$ws.Range("E30").Value = '  +0.08%  '
$ws.Range("E31").Value = '  +2.11%  '
$ws.Range("E32").Value = '  -1.36%  '
$ws.Range("D33").Value = '''3.95'
$ws.Range("E33").Value = '  -1.51%  '
$ws.Range("E34").Value = '  +12.33%  '
$ws.Range("E35").Value = '  +1.15%  '
$ws.Range("E36").Value = '  +1.57%  '
$ws.Range("D37").Value = '''92.70'
$ws.Range("E37").Value = '  -1.25%  '
$ws.Range("E38").Value = '  +7.03%  '
$ws.Range("D39").Value = '1.338.10'
$ws.Range("E39").Value = '  +0.90%  '
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("E41").Value = '  -1.50%  '
$ws.Range("D42").Value = '''14.54'
$ws.Range("E42").Value = '  -3.09%  '
$ws.Range("D43").Value = '''2.27'
$ws.Range("E43").Value = '  -3.52%  '
$ws.Range("D45").Value = '''2.75'
$ws.Range("E45").Value = '  -0.90%  '
$ws.Range("D46").Value = '''0.0521'
$ws.Range("E46").Value = '  +2.30%  '
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("D48").Value = '1.998.53'
$ws.Range("E48").Value = '  -0.29%  '
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("D50").Value = '''0.0665'
$ws.Range("E50").Value = '  +3.72%  '
$ws.Range("D51").Value = '''97.31'
$ws.Range("E51").Value = '  -3.44%  '
